$wb = $excel.ActiveWorkbook

# Update "想去人数" (interest count) values that were refreshed in the data export.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5279
    $ws.Range("F4").Value = 918
}
